$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the redirect "from"/"to" values to be fully-qualified URLs
$ws.Range("A2").Value = "https://example.com/from"
$ws.Range("B2").Value = "https://example.com/to"

# Turn those cells into real hyperlinks pointing at the same URL
$ws.Hyperlinks.Add($ws.Range("A2"), "https://example.com/from")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://example.com/to")

# Move the active selection to A11, matching the saved view state
$ws.Range("A11").Select() | Out-Null
